$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# VERSION 2.0: NEW FIRMWARE UPDATED
# Replace the "core_cm3.h" related error rows (old rows 41-48) and the
# trailing usart rows (old rows 49-61) with the new BOOTLOADER / APP_OLD /
# APP_CURRENT / ESP32_to_STM32 derived error rows.
# ---------------------------------------------------------------------------

# Drop every data row from 41 through the old last row (61); this removes the
# rows completely (not just their content) so re-creating them afterwards
# does not leave stale/phantom cells behind.
$ws.Rows("41:61").Delete()

# Use row 40 (still formatted with style index 2 / 41.4pt wrapped row) as the
# template for formatting the freshly (re)created rows below. Copy only as
# many source columns as will be pasted each time, since pasting always
# stamps the whole shape of the copied range onto the destination.
function Set-Row6 {
    param($r, $a, $b, $c, $d, $e, $f)
    $ws.Range("A40:F40").Copy() | Out-Null
    $dst = $ws.Range("A" + $r + ":F" + $r)
    $dst.PasteSpecial(-4122) | Out-Null
    $dst.RowHeight = 41.4
    $ws.Cells.Item($r,1).Value = $a
    $ws.Cells.Item($r,2).Value = $b
    $ws.Cells.Item($r,3).Value = $c
    $ws.Cells.Item($r,4).Value = $d
    $ws.Cells.Item($r,5).Value = $e
    $ws.Cells.Item($r,6).Value = $f
}

function Set-Row5 {
    param($r, $a, $b, $c, $d, $e)
    $ws.Range("A40:E40").Copy() | Out-Null
    $dst = $ws.Range("A" + $r + ":E" + $r)
    $dst.PasteSpecial(-4122) | Out-Null
    $dst.RowHeight = 41.4
    $ws.Cells.Item($r,1).Value = $a
    $ws.Cells.Item($r,2).Value = $b
    $ws.Cells.Item($r,3).Value = $c
    $ws.Cells.Item($r,4).Value = $d
    $ws.Cells.Item($r,5).Value = $e
}

Set-Row6 41 "stm32f103xx_usart_driver.c" 69  "warning: C23 extension"            "label followed by a declaration is a C23 extension"    "C23-specific syntax" "No action (compiler flag adjustment)"
Set-Row6 42 "stm32f103xx_usart_driver.c" 22  "warning: expression result unused" "expression result unused [-Wunused-value]"              "Unused increment result" "No action (driver file)"
Set-Row6 43 "stm32f103xx_usart_driver.c" 155 "warning: C23 extension"            "label at end of compound statement is a C23 extension"  "C23-specific syntax" "No action (driver file)"
Set-Row6 44 "stm32f103xx_usart_driver.c" 171 "warning: C23 extension"            "label at end of compound statement is a C23 extension"  "C23-specific syntax" "No action (driver file)"
Set-Row6 45 "stm32f103xx_usart_driver.c" 201 "warning: C23 extension"            "label at end of compound statement is a C23 extension"  "C23-specific syntax" "No action (driver file)"
Set-Row6 46 "stm32f103xx_usart_driver.c" 270 "warning: expression result unused" "expression result unused [-Wunused-value]"              "Unused increment result" "No action (driver file)"

Set-Row5 47 "main.c"       331 "error: undeclared function" "call to undeclared function 'USART_Enable'" "Missing function declaration in driver"
Set-Row5 48 "bootloader.c" 29  "error: undeclared function" "call to undeclared function '__disable_irq'" "Missing CMSIS function, avoided by design"
Set-Row5 49 "bootloader.c" 30  "error: undeclared function" "call to undeclared function '__set_MSP'" "Missing CMSIS function, replaced with asm"
Set-Row5 50 "bootloader.c" 32  "error: undeclared function" "call to undeclared function '__enable_irq'" "Missing CMSIS function, avoided by design"
Set-Row5 51 "syscalls.c"   172 "warning: declaration visibility" "declaration of 'struct tms' will not be visible" "Stubbed struct visibility"

# Title cell (A1) keeps the same text; nothing to change there.

# Update the view: scroll so row 39 is near the top and select G43, matching
# the author's on-save cursor/viewport position.
$ws.Range("G43").Select()
$excel.ActiveWindow.ScrollRow = 39
$excel.ActiveWindow.ScrollColumn = 1
